$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45171 = 2023-09-02)
# that needs to be bumped by one day (45172 = 2023-09-03) for every data
# row (rows 2 through 423).
$lastRow = 423
$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45172
